$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule "R40" (cell B11 on the "Rules" sheet) is renamed to "1".
# Force text so the value is stored as a string (shared string),
# matching the original cell's data type, instead of being
# auto-coerced to a number.
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
